# Auto-generated script applying numeric updates to the Leve profit-tracking sheets.
# Each sheet (crafting class) has per-row columns H..N holding Universalis market-price
# derived figures (current average prices, leve sale prices, and computed profits).
# This script refreshes those cached values to the latest scrape, and clears two cells
# whose HQ profit figure is no longer applicable (row now has no HQ price).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2621
$ws.Range("J2").Value = 3026
$ws.Range("L2").Value = 3026
$ws.Range("N2").Value = -3252
$ws.Range("H11").Value = 572.619
$ws.Range("I11").Value = 572.619
$ws.Range("K11").Value = 572.619
$ws.Range("M11").Value = -432.619
$ws.Range("H12").Value = 375.46667
$ws.Range("I12").Value = 368.6154
$ws.Range("J12").Value = 420
$ws.Range("K12").Value = 368.6154
$ws.Range("L12").Value = 420
$ws.Range("M12").Value = -198.6154
$ws.Range("N12").Value = -760
$ws.Range("H18").Value = 427.4
$ws.Range("I18").Value = 427.4
$ws.Range("K18").Value = 427.4
$ws.Range("M18").Value = -143.4
$ws.Range("H42").Value = 901.6
$ws.Range("I42").Value = 54
$ws.Range("J42").Value = 1466.6666
$ws.Range("K42").Value = 162
$ws.Range("L42").Value = 4399.9998
$ws.Range("M42").Value = 68
$ws.Range("N42").Value = -4859.9998
$ws.Range("H69").Value = 18000
$ws.Range("I69").Value = 4000
$ws.Range("K69").Value = 12000
$ws.Range("M69").Value = -11126
$ws.Range("H72").Value = 18000
$ws.Range("I72").Value = 4000
$ws.Range("K72").Value = 36000
$ws.Range("M72").Value = -31632
$ws.Range("H98").Value = 4300.92
$ws.Range("I98").Value = 4016.5
$ws.Range("J98").Value = 5438.6
$ws.Range("K98").Value = 4016.5
$ws.Range("L98").Value = 5438.6
$ws.Range("M98").Value = -2518.5
$ws.Range("N98").Value = -8434.6
$ws.Range("H122").Value = 4300.92
$ws.Range("I122").Value = 4016.5
$ws.Range("J122").Value = 5438.6
$ws.Range("K122").Value = 12049.5
$ws.Range("L122").Value = 16315.8
$ws.Range("M122").Value = -9599.5
$ws.Range("N122").Value = -21215.8
$ws.Range("H125").Value = 3867.5
$ws.Range("I125").Value = 3489
$ws.Range("J125").Value = 4498.3335
$ws.Range("K125").Value = 31401
$ws.Range("L125").Value = 40485.0015
$ws.Range("M125").Value = -28941
$ws.Range("N125").Value = -45405.0015
$ws.Range("H132").Value = 501834.12
$ws.Range("I132").Value = 647484.4
$ws.Range("K132").Value = 1942453.2
$ws.Range("M132").Value = -1939923.2
$ws.Range("H135").Value = 15059.9
$ws.Range("I135").Value = 5149.75
$ws.Range("K135").Value = 46347.75
$ws.Range("M135").Value = -43812.75
$ws.Range("H138").Value = 3790
$ws.Range("J138").Value = 4469.6855
$ws.Range("L138").Value = 13409.0565
$ws.Range("N138").Value = -23689.0565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5367.763
$ws.Range("I61").Value = 7876.636
$ws.Range("K61").Value = 7876.636
$ws.Range("M61").Value = -7664.636
$ws.Range("H74").Value = 3006.1765
$ws.Range("I74").Value = 2441.5625
$ws.Range("J74").Value = 3508.0557
$ws.Range("K74").Value = 2441.5625
$ws.Range("L74").Value = 3508.0557
$ws.Range("M74").Value = -1567.5625
$ws.Range("N74").Value = -5256.0557
$ws.Range("H77").Value = 3006.1765
$ws.Range("I77").Value = 2441.5625
$ws.Range("J77").Value = 3508.0557
$ws.Range("K77").Value = 12207.8125
$ws.Range("L77").Value = 17540.2785
$ws.Range("M77").Value = -7839.8125
$ws.Range("N77").Value = -26276.2785
$ws.Range("H122").Value = 2974.0667
$ws.Range("I122").Value = 2495.3928
$ws.Range("K122").Value = 7486.178400000001
$ws.Range("M122").Value = -5036.178400000001
$ws.Range("H132").Value = 493833.6
$ws.Range("I132").Value = 610120.4
$ws.Range("K132").Value = 1830361.2
$ws.Range("M132").Value = -1827831.2
$ws.Range("H136").Value = 5367.763
$ws.Range("I136").Value = 7876.636
$ws.Range("K136").Value = 23629.908
$ws.Range("M136").Value = -21079.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2345.4
$ws.Range("I54").Value = 2345.4
$ws.Range("K54").Value = 2345.4
$ws.Range("M54").Value = -1861.4
$ws.Range("H86").Value = 5127.5713
$ws.Range("I86").Value = 1176.7778
$ws.Range("K86").Value = 1176.7778
$ws.Range("M86").Value = -53.77780000000007
$ws.Range("H89").Value = 5127.5713
$ws.Range("I89").Value = 1176.7778
$ws.Range("K89").Value = 5883.889
$ws.Range("M89").Value = -267.8890000000001
$ws.Range("H134").Value = 1332750.9
$ws.Range("I134").Value = 1543430
$ws.Range("K134").Value = 4630290
$ws.Range("M134").Value = -4627755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 347.07144
$ws.Range("I7").Value = 92.59090999999999
$ws.Range("K7").Value = 92.59090999999999
$ws.Range("M7").Value = 20.40909000000001
$ws.Range("H99").Value = 18182816
$ws.Range("I99").Value = 30303692
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 30303692
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -30302194
$ws.Range("N99").Value = -4496
$ws.Range("H126").Value = 18182816
$ws.Range("I126").Value = 30303692
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 90911076
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -90908606
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 20589.023
$ws.Range("I132").Value = 26504.688
$ws.Range("K132").Value = 79514.064
$ws.Range("M132").Value = -76984.064
$ws.Range("H133").Value = 74293
$ws.Range("J133").Value = 74290
$ws.Range("L133").Value = 74290
$ws.Range("N133").Value = -79350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 92.703705
$ws.Range("I2").Value = 113.21429
$ws.Range("J2").Value = 70.61539
$ws.Range("K2").Value = 679.28574
$ws.Range("L2").Value = 423.6923400000001
$ws.Range("M2").Value = -566.28574
$ws.Range("N2").Value = -649.6923400000001
$ws.Range("H33").Value = 57.666668
$ws.Range("I33").Value = 31.153847
$ws.Range("K33").Value = 186.923082
$ws.Range("M33").Value = 96.07691800000001
$ws.Range("H98").Value = 6867.1665
$ws.Range("I98").Value = 5067.6665
$ws.Range("K98").Value = 15202.9995
$ws.Range("M98").Value = -13704.9995
$ws.Range("H130").Value = 8007.5
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 2490.5715
$ws.Range("I132").Value = 1024.75
$ws.Range("J132").Value = 2835.4707
$ws.Range("K132").Value = 9222.75
$ws.Range("L132").Value = 25519.2363
$ws.Range("M132").Value = -6692.75
$ws.Range("N132").Value = -30579.2363

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11469.833
$ws.Range("I122").Value = 12985
$ws.Range("K122").Value = 38955
$ws.Range("M122").Value = -36505
$ws.Range("H126").Value = 8628053
$ws.Range("I126").Value = 20836026
$ws.Range("J126").Value = 10659.117
$ws.Range("K126").Value = 62508078
$ws.Range("L126").Value = 31977.351
$ws.Range("M126").Value = -62505608
$ws.Range("N126").Value = -36917.351
$ws.Range("H132").Value = 15154294
$ws.Range("I132").Value = 19610190
$ws.Range("J132").Value = 4244.3335
$ws.Range("K132").Value = 58830570
$ws.Range("L132").Value = 12733.0005
$ws.Range("M132").Value = -58828040
$ws.Range("N132").Value = -17793.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2258.4285
$ws.Range("I22").Value = 1193
$ws.Range("J22").Value = 3057.5
$ws.Range("K22").Value = 1193
$ws.Range("L22").Value = 3057.5
$ws.Range("M22").Value = -898
$ws.Range("N22").Value = -3647.5
$ws.Range("H27").Value = 2258.4285
$ws.Range("I27").Value = 1193
$ws.Range("J27").Value = 3057.5
$ws.Range("K27").Value = 1193
$ws.Range("L27").Value = 3057.5
$ws.Range("M27").Value = -1086
$ws.Range("N27").Value = -3271.5
$ws.Range("H46").Value = 17858632
$ws.Range("J46").Value = 71430930
$ws.Range("L46").Value = 71430930
$ws.Range("N46").Value = -71431306
$ws.Range("H132").Value = 5028
$ws.Range("I132").Value = 5028
$ws.Range("K132").Value = 15084
$ws.Range("M132").Value = -12554
$ws.Range("H139").Value = 131111
$ws.Range("J139").Value = 131111
$ws.Range("L139").Value = 131111
$ws.Range("N139").Value = -141391

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("K81").Value = 2000
$ws.Range("M81").Value = -939
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("K84").Value = 10000
$ws.Range("M84").Value = -4696
$ws.Range("H132").Value = 5931.121
$ws.Range("J132").Value = 10144
$ws.Range("L132").Value = 30432
$ws.Range("N132").Value = -35492
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
